# Correct excel data order mismatch: update a handful of computed values
# in column D ("Free Trues") and column E ("Free Lies") on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value  = 0.2669761724656358
$ws.Range("D3").Value  = 0.2742256193751548
$ws.Range("D4").Value  = 0.07986629951164044
$ws.Range("D5").Value  = 0.1256166703958743
$ws.Range("E6").Value  = 0.2844508345524505
$ws.Range("E7").Value  = 0.1226409328994014
$ws.Range("D8").Value  = -0.001137573191452787
$ws.Range("D9").Value  = 0.0266980013036533
$ws.Range("E10").Value = 0.1897646632683353
$ws.Range("D11").Value = 0.05823102669627696
$ws.Range("E12").Value = 0.2137828053710756
$ws.Range("D13").Value = 0.03691155983262792
$ws.Range("D14").Value = 0.09125660324123014
$ws.Range("D15").Value = 0.1439017924661482
$ws.Range("D16").Value = -0.02335708567310092
$ws.Range("D17").Value = 0.1098382484924759
